$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.619.06"
$ws.Range("E2").Value = "  +2.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.97"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.29"
$ws.Range("E5").Value = "  +1.04%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  +0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4909"
$ws.Range("E7").Value = "  -0.21%  "

# Row 8
$ws.Range("E8").Value = "  +0.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06774"
$ws.Range("E9").Value = "  +2.22%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.888.41"
$ws.Range("E10").Value = "  +0.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.26"
$ws.Range("E11").Value = "  +3.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07242"
$ws.Range("E12").Value = "  +0.94%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.02"
$ws.Range("E13").Value = "  +4.65%  "

# Row 14
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6769"
$ws.Range("E14").Value = "  +1.35%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.050"
$ws.Range("E15").Value = "  +3.30%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.606.99"
$ws.Range("E16").Value = "  +2.16%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007963"
$ws.Range("E17").Value = "  +1.75%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.14"
$ws.Range("E19").Value = "  +2.71%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.130.25"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("E21").Value = "  +0.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.815"
$ws.Range("E22").Value = "  +0.69%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "188.46"
$ws.Range("E23").Value = "  +33.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.075"
$ws.Range("E24").Value = "  +3.70%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.331"
$ws.Range("E25").Value = "  +2.42%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.80"
$ws.Range("E26").Value = "  +3.65%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").Value = "  +12.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.905"
$ws.Range("E28").Value = "  -0.19%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.399"
$ws.Range("E29").Value = "  +0.69%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.330"
$ws.Range("E30").Value = "  +2.90%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09047"
$ws.Range("E31").Value = "  +3.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.010"
$ws.Range("E32").Value = "  +0.58%  "

# Row 33
$ws.Range("E33").Value = "  +3.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7521"
$ws.Range("E34").Value = "  +4.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.111"
$ws.Range("E35").Value = "  -0.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.750"
$ws.Range("E36").Value = "  +3.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01841"
$ws.Range("E37").Value = "  +2.01%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.672"
$ws.Range("E38").Value = "  -0.94%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.146"
$ws.Range("E39").Value = "  -0.59%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9338"
$ws.Range("E40").Value = "  -0.52%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4423"
$ws.Range("E41").Value = "  +4.39%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.71"
$ws.Range("E42").Value = "  +1.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("E44").Value = "  -0.24%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.589"
$ws.Range("E45").Value = "  +3.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1340"
$ws.Range("E46").Value = "  +5.55%  "

# Row 47
$ws.Range("E47").Value = "  +2.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.439"
$ws.Range("E48").Value = "  +7.31%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.690"
$ws.Range("E49").Value = "  +4.87%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3926"
$ws.Range("E50").Value = "  +4.10%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.52"
$ws.Range("E51").Value = "  +2.68%  "
